$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$style = $ws.Range("D2").Style
$ws.Range("D2").Value = "'29.865.66"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  +0.38%  "

# Row 3
$style = $ws.Range("D3").Style
$ws.Range("D3").Value = "'1.893.20"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  +0.35%  "

# Row 4
$style = $ws.Range("D4").Style
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'0.7820"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -1.20%  "

# Row 6
$style = $ws.Range("D6").Style
$ws.Range("D6").Value = "'243.60"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  +0.73%  "

# Row 7
$style = $ws.Range("D7").Style
$ws.Range("D7").Value = "'1.001"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$style = $ws.Range("D8").Style
$ws.Range("D8").Value = "'0.3136"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  -0.92%  "

# Row 9
$style = $ws.Range("D9").Style
$ws.Range("D9").Value = "'25.71"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  +1.43%  "

# Row 10
$style = $ws.Range("D10").Style
$ws.Range("D10").Value = "'0.07262"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  +3.75%  "

# Row 11
$style = $ws.Range("D11").Style
$ws.Range("D11").Value = "'0.08094"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  +0.68%  "

# Row 12
$style = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.7732"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  +0.99%  "

# Row 13
$style = $ws.Range("D13").Style
$ws.Range("D13").Value = "'5.466"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  +3.20%  "

# Row 14
$style = $ws.Range("D14").Style
$ws.Range("D14").Value = "'1.870.51"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -0.62%  "

# Row 15
$style = $ws.Range("D15").Style
$ws.Range("D15").Value = "'94.00"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  +2.28%  "

# Row 16
$style = $ws.Range("D16").Style
$ws.Range("D16").Value = "'6.207"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  +5.04%  "

# Row 17
$style = $ws.Range("D17").Style
$ws.Range("D17").Value = "'29.879.56"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  +0.43%  "

# Row 18
$style = $ws.Range("D18").Style
$ws.Range("D18").Value = "'13.92"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  +0.47%  "

# Row 19
$style = $ws.Range("D19").Style
$ws.Range("D19").Value = "'246.60"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  +1.40%  "

# Row 20
$style = $ws.Range("D20").Style
$ws.Range("D20").Value = "'0.000007800"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  +1.51%  "

# Row 21
$style = $ws.Range("D21").Style
$ws.Range("D21").Value = "'2.159.17"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +1.31%  "

# Row 22
$ws.Range("E22").Value = "  +0.04%  "

# Row 23
$style = $ws.Range("D23").Style
$ws.Range("D23").Value = "'8.131"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -0.73%  "

# Row 24
$style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'1.002"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  +0.10%  "

# Row 25
$style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'0.1595"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -5.32%  "

# Row 26
$style = $ws.Range("D26").Style
$ws.Range("D26").Value = "'9.437"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  +1.50%  "

# Row 27
$style = $ws.Range("D27").Style
$ws.Range("D27").Value = "'164.44"
$ws.Range("D27").Style = $style

# Row 28
$style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'18.74"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  +0.60%  "

# Row 29
$style = $ws.Range("D29").Style
$ws.Range("D29").Value = "'2.024"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -1.14%  "

# Row 30
$ws.Range("E30").Value = "  +3.36%  "

# Row 31
$style = $ws.Range("D31").Style
$ws.Range("D31").Value = "'1.543"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  +0.61%  "

# Row 32
$style = $ws.Range("D32").Style
$ws.Range("D32").Value = "'4.475"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  +2.07%  "

# Row 33
$style = $ws.Range("D33").Style
$ws.Range("D33").Value = "'0.05568"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -1.88%  "

# Row 34
$style = $ws.Range("D34").Style
$ws.Range("D34").Value = "'4.069"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  +0.58%  "

# Row 35
$style = $ws.Range("D35").Style
$ws.Range("D35").Value = "'1.240"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -1.53%  "

# Row 36
$style = $ws.Range("D36").Style
$ws.Range("D36").Value = "'0.7537"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  +2.46%  "

# Row 37
$style = $ws.Range("D37").Style
$ws.Range("D37").Value = "'1.002"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  +0.15%  "

# Row 38
$style = $ws.Range("D38").Style
$ws.Range("D38").Value = "'2.683"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  +1.85%  "

# Row 39
$style = $ws.Range("D39").Style
$ws.Range("D39").Value = "'0.01933"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  +1.34%  "

# Row 40
$style = $ws.Range("D40").Style
$ws.Range("D40").Value = "'2.798"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  +1.12%  "

# Row 41
$style = $ws.Range("D41").Style
$ws.Range("D41").Value = "'1.139.49"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  +11.48%  "

# Row 42
$style = $ws.Range("D42").Style
$ws.Range("D42").Value = "'0.4462"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  +1.56%  "

# Row 43
$style = $ws.Range("D43").Style
$ws.Range("D43").Value = "'73.97"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  +2.28%  "

# Row 44
$style = $ws.Range("D44").Style
$ws.Range("D44").Value = "'5.971"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  +2.33%  "

# Row 45
$style = $ws.Range("D45").Style
$ws.Range("D45").Value = "'0.8523"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +1.98%  "

# Row 46
$style = $ws.Range("D46").Style
$ws.Range("D46").Value = "'1.001"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  +0.07%  "

# Row 47
$style = $ws.Range("D47").Style
$ws.Range("D47").Value = "'1.888"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  +1.90%  "

# Row 48
$style = $ws.Range("D48").Style
$ws.Range("D48").Value = "'3.134"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +8.22%  "

# Row 49
$style = $ws.Range("D49").Style
$ws.Range("D49").Value = "'102.40"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  +0.04%  "

# Row 50
$style = $ws.Range("D50").Style
$ws.Range("D50").Value = "'7.529"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  +1.58%  "

# Row 51
$style = $ws.Range("D51").Style
$ws.Range("D51").Value = "'9.745"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  -1.29%  "
